$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 517, shifting rows 517:629 down to 518:630
$ws.Rows("517:517").Insert()

# Populate the newly inserted row 517 with the new data record
$ws.Range("A517").Value = 9
$ws.Range("B517").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C517").Value = "Metropolitana"
$ws.Range("D517").Value = 45211
$ws.Range("E517").Value = 13
$ws.Range("F517").Value = 100112044
$ws.Range("G517").Value = "Perejil"
$ws.Range("H517").Value = "Sin especificar"
$ws.Range("I517").Value = "Primera"
$ws.Range("J517").Value = 70
$ws.Range("K517").Value = 14000
$ws.Range("L517").Value = 15000
$ws.Range("M517").Value = 14500
$ws.Range("N517").Value = "`$/docena de atados"
$ws.Range("O517").Value = "Región Metropolitana"
$ws.Range("P517").Value = 4833
$ws.Range("Q517").Value = 3
$ws.Range("R517").Value = "Hortaliza"
